# Update MSME definitions: replace literal "<br/>" separators inside the
# "Number of employees" and "Turnover" definition cells (Small/Medium/Large
# rows) with real line breaks (newline characters), matching the XML diff.
#
# Several of these strings start with "=" (e.g. "=<19 Manufacturing"), so
# assigning them straight to .Value/.Formula would make Excel parse them as
# formulas (and they'd error out, since "<19" etc. isn't valid formula
# syntax). Instead we build each string as a CHAR(10)-joined text formula,
# then convert it in place to a literal value with Copy + PasteSpecial
# (values only) - exactly like typing the formula in and pasting-as-values,
# so the result is plain text, not a formula, and no extra number
# formatting gets attached to the cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-LiteralTextWithBreaks {
    param(
        [string]$CellRef,
        [string[]]$Lines
    )

    $quoted = $Lines | ForEach-Object { '"' + $_ + '"' }
    $formula = "=" + ($quoted -join "&CHAR(10)&")

    $range = $ws.Range($CellRef)
    $range.Formula = $formula
    $range.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
}

# Row 24 - Small
Set-LiteralTextWithBreaks "B24" @("=<19 Manufacturing", "=<49 Service")
Set-LiteralTextWithBreaks "D24" @("=< MNT 250 Millionlion Manufacturing", "=< MNT 1 Billionlion Service")

# Row 25 - Medium
Set-LiteralTextWithBreaks "B25" @("=<149 Wholesale trade", "=<199 Retail Trade", "=<199 Manufacturing")
Set-LiteralTextWithBreaks "D25" @("=< MNT 1.5 Billionlion Wholesale trade", "=< MNT 1.5 Billionlion Retail Trade", "=< MNT 1.5 Billionlion Manufacturing")

# Row 26 - Large
Set-LiteralTextWithBreaks "B26" @(">149 Wholesale trade", ">199 Retail Trade", ">199 Manufacturing")
Set-LiteralTextWithBreaks "D26" @("> MNT 1.5 Billionlion Wholesale trade", "> MNT 1.5 Billionlion Retail Trade", "> MNT 1.5 Billionlion Manufacturing")
